# Fix bug in update_working_list: two rows had wrong dates (13.12./14.12.
# duplicated from an earlier entry) -- they should read 23.12./24.12. Also
# add the missing 25.12.2024 work log entries and fix the running total.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the mis-dated rows (were "13.12.2024"/"14.12.2024", duplicates of
#     row 43/44's dates -- should be "23.12.2024"/"24.12.2024") ---
$ws.Range("A46").Value = "23.12.2024"
$ws.Range("A47").Value = "24.12.2024"

# --- New entries for 25.12.2024 ---
$ws.Range("A50").Value = "25.12.2024"
$ws.Range("B50").Value = 8
$ws.Range("C50").Value = "Verbesserung, Testen und Debuggen des Prototypen"

$ws.Range("B51").Value = 2
$ws.Range("C51").Value = "Vergleich des aktuellen Prototypen mit dem Original"

$wb.Save()
